$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 (Item 13 - 81044/12-24-6 hook-up wire) — copy row 13's
# formatting (borders, fonts, number formats) down to row 14 first so the
# new row matches the look of the existing BOM rows.
$ws.Range("A13:H13").Copy($ws.Range("A14:H14"))

# Fill in the new row's data.
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "`t81044/12-24-6-DS-ND"
$ws.Range("D14").Value = "`t81044/12-24-6"
$ws.Range("E14").Value = '81044/12-24-6 FEET HOOK-UP WIRE'
$ws.Range("F14").Value = 0.98
$ws.Range("G14").Value = 0.98

# Add the Digikey link for the new part, then restore the row-13 look for
# column H (Hyperlinks.Add applies Excel's default Hyperlink style, which
# would otherwise strip the table border) and set the visible text back to
# the link target (matching how it was originally inserted).
$ws.Hyperlinks.Add($ws.Range("H14"), 'https://www.digikey.ca/en/products/detail/te-connectivity-aerospace,-defense-and-marine/81044%2F12-24-6/6071077')
$ws.Range("H13").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = 'https://www.digikey.ca/en/products/detail/te-connectivity-aerospace,-defense-and-marine/81044%2F12-24-6/6071077'

# Extend the total-price formula to include the new row.
$ws.Range("G16").Formula = "=F2*B2+F3*B3+F4*B4+F5*B5+F6*B6+F7*B7+F8*B8+B9*F9+B10*F10+B11*F11+B12*F12+B13*F13+B14*F14"

# Update the footnote to mention both optional items.
$ws.Range("C17").Value = '*Items 12 & 13 are optional '

# Match the author's final selection.
$ws.Range("C17").Select()
